$q3Data = @(
    @(0, '519674', '银河创新成长混合A', '142.10', '94.51', '7.47', '10.6149', 7),
    @(1, '161903', '万家行业优选混合（LOF）', '83.48', '93.78', '7.36', '6.1441', 4),
    @(2, '501207', '华夏创新未来混合（LOF）', '44.72', '91.75', '4.09', '1.8290', 9),
    @(3, '008120', '万家自主创新混合A', '22.86', '93.71', '7.43', '1.6985', 5),
    @(4, '014143', '银河创新成长混合C', '20.25', '94.51', '7.47', '1.5127', 7),
    @(5, '506005', '博时科创板三年定期开放混合', '21.62', '99.12', '5.98', '1.2929', 3),
    @(6, '000031', '华夏复兴混合A', '24.25', '88.23', '4.79', '1.1616', 8),
    @(7, '290011', '泰信中小盘精选混合', '11.26', '94.48', '9.11', '1.0258', 8),
    @(8, '007490', '南方信息创新混合A', '16.36', '91.57', '5.41', '0.8851', 3),
    @(9, '005311', '万家经济新动能混合A', '10.59', '93.55', '7.19', '0.7614', 6),
    @(10, '007349', '华夏科技创新混合A', '11.68', '90.65', '4.06', '0.4742', 10),
    @(11, '005312', '万家经济新动能混合C', '5.65', '93.55', '7.19', '0.4062', 6),
    @(12, '013840', '银华集成电路混合A', '7.35', '92.35', '4.22', '0.3102', 7),
    @(13, '013962', '华夏创新视野一年持有混合A', '7.47', '88.70', '4.04', '0.3018', 9),
    @(14, '013963', '华夏创新视野一年持有混合C', '6.46', '88.70', '4.04', '0.2610', 9),
    @(15, '010106', '华夏核心科技6个月定期开放混合A', '6.17', '81.38', '4.04', '0.2493', 5),
    @(16, '010518', '华夏先锋科技一年定期开放混合A', '5.64', '89.54', '4.06', '0.2290', 10),
    @(17, '013841', '银华集成电路混合C', '5.23', '92.35', '4.22', '0.2207', 7),
    @(18, '008418', '惠升惠泽灵活配置混合A', '7.69', '77.30', '2.73', '0.2099', 8),
    @(19, '540010', '汇丰晋信科技先锋股票', '3.39', '93.51', '6.09', '0.2065', 7),
    @(20, '012650', '博时半导体主题混合A', '4.43', '93.72', '4.44', '0.1967', 7),
    @(21, '008633', '万家科技创新混合A', '2.44', '93.37', '7.30', '0.1781', 7),
    @(22, '008121', '万家自主创新混合C', '2.24', '93.71', '7.43', '0.1664', 5),
    @(23, '012651', '博时半导体主题混合C', '3.30', '93.72', '4.44', '0.1465', 7),
    @(24, '008716', '鹏华优质回报两年定期开放混合', '5.93', '66.94', '2.37', '0.1405', 10),
    @(25, '007491', '南方信息创新混合C', '2.48', '91.57', '5.41', '0.1342', 3),
    @(26, '008634', '万家科技创新混合C', '1.71', '93.37', '7.30', '0.1248', 7),
    @(27, '001188', '鹏华改革红利股票', '2.63', '88.79', '4.65', '0.1223', 3),
    @(28, '588100', '嘉实上证科创板新一代信息技术ETF', '3.66', '99.40', '2.76', '0.1010', 10),
    @(29, '002580', '泰信鑫选灵活配置混合C', '1.06', '94.19', '9.44', '0.1001', 8),
    @(30, '006366', '兴业安保优选混合', '1.26', '87.20', '7.31', '0.0921', 2),
    @(31, '360011', '光大保德信动态优选混合', '1.64', '74.29', '5.52', '0.0905', 6),
    @(32, '007854', '光大保德信景气先锋混合', '1.50', '73.83', '5.89', '0.0884', 5),
    @(33, '001970', '泰信鑫选灵活配置混合A', '0.79', '94.19', '9.44', '0.0746', 8),
    @(34, '010519', '华夏先锋科技一年定期开放混合C', '1.78', '89.54', '4.06', '0.0723', 10),
    @(35, '000522', '华润元大信息传媒科技混合', '1.42', '71.23', '3.68', '0.0523', 9),
    @(36, '588260', '华安上证科创板新一代信息技术ETF', '1.73', '98.52', '2.75', '0.0476', 10),
    @(37, '011077', '汇丰晋信创新先锋股票', '0.76', '94.21', '5.97', '0.0454', 8),
    @(38, '010107', '华夏核心科技6个月定期开放混合C', '1.09', '81.38', '4.04', '0.0440', 5),
    @(39, '007350', '华夏科技创新混合C', '0.88', '90.65', '4.06', '0.0357', 10),
    @(40, '002863', '金信深圳成长灵活配置混合', '0.61', '93.95', '4.16', '0.0254', 10),
    @(41, '008419', '惠升惠泽灵活配置混合C', '0.76', '77.30', '2.73', '0.0207', 8),
    @(42, '002772', '光大保德信产业新动力灵活配置混合', '0.25', '87.41', '7.55', '0.0189', 4),
    @(43, '001574', '中海混改红利主题精选灵活配置混合', '0.24', '92.81', '6.31', '0.0151', 6),
    @(44, '015073', '华夏复兴混合C', '0.04', '88.23', '4.79', '0.0019', 8)
)

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for "2022-Q3" and shift the
#    existing quarters down by one row, renumbering the index column.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows(2).Insert()

# The inserted row doesn't inherit formatting from the row below it,
# so copy column-A's style (bold + border) from the row that used to
# be row 2 (now row 3) onto the new row 2.
$wsTotal.Range("A3").Copy($wsTotal.Range("A2"))
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 45
$wsTotal.Range("D2").Value = 31.93

# Renumber the 0-based index column for the rows that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
}

# ------------------------------------------------------------------
# 2) Add a new "2022-Q3" sheet (copied from "2022-Q2" so it inherits
#    the same layout/formatting), placed right before "2022-Q2", then
#    overwrite its data with the 2022-Q3 fund holdings.
# ------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

foreach ($row in $q3Data) {
    $r = $row[0] + 2
    $wsQ3.Cells.Item($r, 1).Value = $row[0]
    $wsQ3.Cells.Item($r, 2).Value = "'" + $row[1]
    $wsQ3.Cells.Item($r, 3).Value = $row[2]
    $wsQ3.Cells.Item($r, 4).Value = "'" + $row[3]
    $wsQ3.Cells.Item($r, 5).Value = "'" + $row[4]
    $wsQ3.Cells.Item($r, 6).Value = "'" + $row[5]
    $wsQ3.Cells.Item($r, 7).Value = "'" + $row[6]
    $wsQ3.Cells.Item($r, 8).Value = $row[7]
}

# The copied sheet had 85 data rows (1 header + 84); the 2022-Q3 sheet
# only has 46 (1 header + 45), so clear the leftover rows.
$wsQ3.Range("A47:H85").Clear()

Write-Host "done"

